# Applies the "redescente LRM hub" update to the ValueSet workbook:
#  - Sets the "Experimental" flag value (row 7, column B on the Metadata
#    sheet) to the text "false"
#  - Refreshes the "Date" value (row 8, column B on the Metadata sheet)
#    to the new publication timestamp

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Experimental -> "false" (must land as literal text, not a Boolean, so
# build it via a formula and then flatten the formula to a static value
# with a copy/paste-special so the cell keeps its original style/type).
$expCell = $ws.Range("B7")
$expCell.Formula = "=""false"""
$expCell.Copy()
$expCell.PasteSpecial(-4163)

# Date -> new timestamp
$ws.Range("B8").Value = "2025-11-04T10:04:56+00:00"
